$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated price/volume cells remain stored as text (matching the
# original inline-string cells) rather than being auto-converted to numbers
# or percentages by Excel.
$cells = @(
    @('D2', '330.63'),
    @('E2', '-0.42%'),
    @('D3', '41.60'),
    @('E3', '0.75%'),
    @('D4', '5.689'),
    @('E4', '0.19%'),
    @('D5', '0.08352'),
    @('E5', '3.40%'),
    @('D6', '8.808'),
    @('E6', '0.58%'),
    @('D7', '2.005'),
    @('E7', '-1.95%'),
    @('D8', '4.467'),
    @('E8', '-1.66%'),
    @('D10', '0.9243'),
    @('E10', '0.09%'),
    @('D11', '0.1288'),
    @('E11', '2.53%'),
    @('D12', '0.1971'),
    @('E12', '1.31%'),
    @('D13', '0.09469'),
    @('E13', '1.37%'),
    @('D14', '0.03840'),
    @('E14', '2.80%'),
    @('E15', '0.92%'),
    @('D16', '0.001305'),
    @('E16', '-0.93%'),
    @('D17', '0.006102'),
    @('E17', '-1.38%'),
    @('D18', '3.428'),
    @('E18', '1.86%'),
    @('E19', '0.65%'),
    @('D20', '8.772'),
    @('E20', '-0.01%'),
    @('D21', '0.1361'),
    @('E21', '-4.06%'),
    @('D22', '0.2507'),
    @('E22', '-5.58%'),
    @('D23', '0.04406'),
    @('E23', '-0.66%'),
    @('D24', '0.001272'),
    @('E24', '0.84%'),
    @('D25', '0.004382'),
    @('E25', '1.15%'),
    @('E26', '-1.86%'),
    @('D27', '0.0003988'),
    @('E27', '-0.13%'),
    @('D39', '0.02840'),
    @('E39', '-1.28%'),
    @('D40', '0.05524'),
    @('E40', '0.61%'),
    @('D41', '0.007965'),
    @('E41', '2.07%'),
    @('D42', '0.1433'),
    @('E42', '0.97%'),
    @('D43', '0.008998'),
    @('E43', '-9.59%'),
    @('D44', '0.002061'),
    @('E44', '-3.53%'),
    @('E45', '5.56%'),
    @('D46', '0.00006929'),
    @('E46', '1.65%'),
    @('E47', '-0.31%'),
    @('D48', '0.003464'),
    @('E48', '14.82%'),
    @('E49', '-0.30%'),
    @('D50', '0.00002099'),
    @('E50', '-0.31%'),
    @('D51', '0.0001999'),
    @('E51', '-0.31%')
)

foreach ($pair in $cells) {
    $addr = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}
